$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell B11 changes from the text "R40" to the text "1". A plain
# Range.Value assignment of "1" would be auto-coerced to the number 1 by
# Excel's normal type inference, which would also change the cell's style
# (text-forcing via NumberFormat/quotePrefix allocates a brand new cellXf).
# Routing the literal through a text-producing formula and pasting only the
# value keeps the result a genuine string (shared string table entry) while
# leaving B11's existing style (s="23") untouched.
$scratch = $ws.Range("Z1")
$scratch.Formula = "=""1"""
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()
